$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy row 19's formatting down into the new row 20 so the new data
# inherits the same column styles (date format, wrap/vertical-centered
# number cells) used throughout the table.
$ws.Range("A19:M19").Copy() | Out-Null
$ws.Range("A20:M20").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# New daily row for 2020-04-04 (day 19 of the series).
$ws.Range("A20").Formula = "=A19+1"
$ws.Range("B20").Formula = "=B19+1"
$ws.Range("C20").Value = 105
$ws.Range("D20").Value = 95
$ws.Range("E20").Value = 1426
$ws.Range("F20").Value = 1626
$ws.Range("G20").Value = "#N/A"
$ws.Range("H20").Value = "#N/A"
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 2
$ws.Range("L20").Value = 0
$ws.Range("K20").Formula = "=K19+L20"
$ws.Range("M20").Value = 54

$ws.Range("N20").Select()
